# Update "F" column (想去人数 / want-to-go counts) figures across the
# four worksheets to match the newly generated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibition) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 169
$ws.Range("F8").Value = 805
$ws.Range("F9").Value = 4201
$ws.Range("F11").Value = 174
$ws.Range("F13").Value = 6074
$ws.Range("F14").Value = 62
$ws.Range("F16").Value = 2333
$ws.Range("F18").Value = 167
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 9141
$ws.Range("F21").Value = 41
$ws.Range("F22").Value = 2455
$ws.Range("F23").Value = 196
$ws.Range("F24").Value = 2310
$ws.Range("F25").Value = 2434
$ws.Range("F26").Value = 1391
$ws.Range("F27").Value = 242
$ws.Range("F28").Value = 1958
$ws.Range("F31").Value = 329
$ws.Range("F33").Value = 38
$ws.Range("F36").Value = 52
$ws.Range("F38").Value = 1222
$ws.Range("F41").Value = 97
$ws.Range("F43").Value = 1530
$ws.Range("F44").Value = 2510
$ws.Range("F45").Value = 923
$ws.Range("F46").Value = 294
$ws.Range("F47").Value = 1250
$ws.Range("F48").Value = 18

# ---- 演出 (Performance) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F21").Value = 29

# ---- 本地生活 (Local Life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 691
$ws.Range("F3").Value = 889

# ---- 全部类型 (All Types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 691
$ws.Range("F5").Value = 889
$ws.Range("F12").Value = 169
$ws.Range("F13").Value = 805
$ws.Range("F14").Value = 4201
$ws.Range("F15").Value = 174
$ws.Range("F16").Value = 6074
$ws.Range("F17").Value = 62
$ws.Range("F19").Value = 2333
$ws.Range("F20").Value = 167
$ws.Range("F21").Value = 472
$ws.Range("F22").Value = 9141
$ws.Range("F24").Value = 2455
$ws.Range("F25").Value = 2310
$ws.Range("F26").Value = 1391
$ws.Range("F27").Value = 242
$ws.Range("F28").Value = 1959
$ws.Range("F31").Value = 329
$ws.Range("F32").Value = 38
$ws.Range("F35").Value = 52
$ws.Range("F39").Value = 97
$ws.Range("F41").Value = 1530
$ws.Range("F42").Value = 2510
$ws.Range("F43").Value = 923
$ws.Range("F44").Value = 294
$ws.Range("F47").Value = 29
$ws.Range("F48").Value = 1250
$ws.Range("F49").Value = 18
